$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows above row 2, shifting the existing B2:D17 block
# down to B4:D19 (matches the diff: data moved from rows 2-17 to 4-19).
$ws.Rows("2:3").Insert()

# Update the selection to match the post-edit state.
$ws.Range("E7").Select()
